# "move to js proto"
#
# - rename the "id" reference column header on dynamic_table_2 (sheet2)
#   from "static_table_1_id" to "static_table_1_refid"
# - flip which sheet/cell is active: previously dynamic_table_2 was the
#   active tab (selection D10); now static_table_1 should be the active
#   tab (selection D18), and dynamic_table_2's own remembered selection
#   moves to D2.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("static_table_1")
$ws2 = $wb.Worksheets.Item("dynamic_table_2")

# Rename the shared-string header text used by D1 on dynamic_table_2.
$ws2.Range("D1").Value = "static_table_1_refid"

# Set dynamic_table_2's remembered selection first (while it's still the
# active sheet) so selecting on it doesn't clobber the final active tab.
$ws2.Range("D2").Select()

# Now make static_table_1 the active sheet/tab and set its selection.
$ws1.Activate()
$ws1.Range("D18").Select()
